# Append 4 new reconciliation rows (164-167) to Sheet1, matching the
# source data that was appended to the underlying report.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 164
$ws.Cells.Item(164, 1).Value = 237674841555
$ws.Cells.Item(164, 2).Value = "BEATRICE TCHAMTIEU EPSE NGAMENI"
$ws.Cells.Item(164, 3).Value = "Rte_5"
$ws.Cells.Item(164, 4).Value = "Hopital General Douala"
$ws.Cells.Item(164, 5).Value = 150379.6
$ws.Cells.Item(164, 6).Value = 790555
$ws.Cells.Item(164, 7).Value = 640175.4
$ws.Cells.Item(164, 8).Value = 5.257062793091616
$ws.Cells.Item(164, 9).Value = "Ndogbong"

# Row 165
$ws.Cells.Item(165, 1).Value = 237674899678
$ws.Cells.Item(165, 2).Value = "VIVIANE MADJUIMEKEM FOMEKONG"
$ws.Cells.Item(165, 3).Value = "Rte_5"
$ws.Cells.Item(165, 4).Value = "Hopital General Douala"
$ws.Cells.Item(165, 5).Value = 159035.175
$ws.Cells.Item(165, 6).Value = 507599
$ws.Cells.Item(165, 7).Value = 348563.825
$ws.Cells.Item(165, 8).Value = 3.191740443584258
$ws.Cells.Item(165, 9).Value = "Ndogbong"

# Row 166
$ws.Cells.Item(166, 1).Value = 237676439452
$ws.Cells.Item(166, 2).Value = "SAGNOU BRINDA JOSELINE _DIGITAL BUSINESS SARL"
$ws.Cells.Item(166, 3).Value = "Rte_5"
$ws.Cells.Item(166, 4).Value = "Hopital General Douala"
$ws.Cells.Item(166, 5).Value = 87377.35000000001
$ws.Cells.Item(166, 6).Value = 9597
$ws.Cells.Item(166, 7).Value = -77780.35000000001
$ws.Cells.Item(166, 8).Value = 0.1098339558249363
$ws.Cells.Item(166, 9).Value = "Ndogbong"

# Row 167
$ws.Cells.Item(167, 1).Value = 237676695935
$ws.Cells.Item(167, 2).Value = "MFS MAKEPE MATURITE"
$ws.Cells.Item(167, 3).Value = "Rte_5"
$ws.Cells.Item(167, 4).Value = "Hopital General Douala"
$ws.Cells.Item(167, 5).Value = 88445
$ws.Cells.Item(167, 6).Value = 312534
$ws.Cells.Item(167, 7).Value = 224089
$ws.Cells.Item(167, 8).Value = 3.53365368307988
$ws.Cells.Item(167, 9).Value = "Ndogbong"
